# Batch-user-signup template: rename the trailing "Unidade(s)*" header to a
# non-mandatory "Unidade(s)" (drops the "required" asterisk), clear the
# leftover bold/"applyFont" style that had been sitting on A1, and leave the
# sheet's active selection on E2 (the cargo/função dropdown cell) the way the
# author left it after testing the batch-add flow.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Unidade(s)*" -> "Unidade(s)" : this is the header cell of the last column
# of the Usuarios table (F1), so editing the cell text also renames the
# ListColumn / shared string in one shot.
$ws.Range("F1").Value = "Unidade(s)"

# A1 had an explicit (but effectively blank) font style applied to it;
# clear it so the cell falls back to the default "Normal" style again.
$ws.Range("A1").ClearFormats()

# Leave the selection on E2 (Cargo/Função column) instead of A2.
$ws.Range("E2").Select()
